$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 306
$ws1.Range("F3").Value = 97
$ws1.Range("F4").Value = 1238
$ws1.Range("F5").Value = 617

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 306
$ws4.Range("F3").Value = 97
$ws4.Range("F4").Value = 1238
$ws4.Range("F6").Value = 617
